$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2-97) forward by 14 days (two weeks)
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $raw = $cell.Value2()
    $cell.Value = $raw + 14
}

# Record the newly observed wind production (MW) for the corresponding hours
$ws.Cells.Item(41, 2).Value = 122
$ws.Cells.Item(42, 2).Value = 120
$ws.Cells.Item(43, 2).Value = 130
